{"js": "// Replace the word \"Resource\" with \"Risk\" in \"Resource assessment \",\n// e.g. \"Resource assessment \" -> \"Risk assessment \".\n//\n// Word leaves the \"_GoBack\" bookmark (last-edit marker) positioned right\n// after the newly typed text, which splits the paragraph's single run into\n// two runs (before/after the bookmark).\n\nconst body = context.document.body;\n\nconst results = body.search(\"Resource\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const found = results.items[0];\n\n  // Stable collapsed range marking where \"Resource\" started (and where the\n  // replacement text will end up), captured before the delete.\n  const insertionPoint = found.getRange(\"Start\");\n\n  found.delete();\n  await context.sync();\n\n  // Move \"_GoBack\" so it marks the point right after the replacement text,\n  // which is what splits the surrounding text into two runs on save.\n  context.document.deleteBookmark(\"_GoBack\");\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n\n  // Type \"Risk\" in place of the deleted \"Resource\".\n  insertionPoint.insertText(\"Risk\", Word.InsertLocation.before);\n  await context.sync();\n}\n", "ps1": "# Replace the word \"Resource\" with \"Risk\" in \"Resource assessment \",\n# e.g. \"Resource assessment \" -> \"Risk assessment \".\n#\n# Word leaves the \"_GoBack\" bookmark (last-edit marker) positioned right\n# after the newly typed text, which splits the paragraph's single run into\n# two runs (before/after the bookmark).\n\n$d = $word.ActiveDocument\n\n# Locate \"Resource\" in the document body.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Resource\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $true\n$rng.Find.MatchWildcards = $false\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n$found = $rng.Find.Execute()\n\nif ($found) {\n    $editPoint = $rng.End\n\n    # Move \"_GoBack\" so it marks the point right after the replacement text,\n    # which is what splits the surrounding text into two runs on save.\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks.Item(\"_GoBack\").Delete()\n    }\n    $d.Bookmarks.Add(\"_GoBack\", $d.Range($editPoint, $editPoint))\n\n    # Replace \"Resource\" with \"Risk\".\n    $rng.Text = \"Risk\"\n}\n"}
